$d = $word.ActiveDocument

# Split the combined "Latex demo / Overleaf" bullet into two separate
# list items: the sentence about the Overleaf report becomes its own
# bullet point (inheriting the same ListParagraph / numPr formatting).
$found = $d.Content.Find.Execute(
    "sections. Ursula set up a group report on Overleaf.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sections.^pUrsula set up a group report on Overleaf.", 2)

if (-not $found) {
    throw "Could not find the target sentence to split into a new bullet."
}
